$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.582.66'
$ws.Range("E2").Value = '  +3.70%  '
$ws.Range("D3").Value = '1.795.98'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.66'
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5310'
$ws.Range("E7").Value = '  -0.79%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3771'
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07531'
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.62'
$ws.Range("E10").Value = '  -1.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.121'
$ws.Range("E11").Value = '  +0.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.18'
$ws.Range("E12").Value = '  +1.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.002'
$ws.Range("E13").Value = '  +0.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.205'
$ws.Range("E14").Value = '  +0.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.494'
$ws.Range("E15").Value = '  +5.82%  '
$ws.Range("D16").Value = '1.793.79'
$ws.Range("E16").Value = '  +0.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '90.41'
$ws.Range("E17").Value = '  -0.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001068'
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06467'
$ws.Range("E19").Value = '  -0.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.33'
$ws.Range("E21").Value = '  +2.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.926'
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("D23").Value = '28.609.38'
$ws.Range("E23").Value = '  +3.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.20'
$ws.Range("E24").Value = '  -0.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.098'
$ws.Range("E25").Value = '  +0.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.99'
$ws.Range("E26").Value = '  +3.59%  '
$ws.Range("E27").Value = '  +0.22%  '
$ws.Range("B28").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C28").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D28").Value = '2.001.14'
$ws.Range("E28").Value = '  +0.36%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.361'
$ws.Range("E29").Value = '  -1.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.93'
$ws.Range("E30").Value = '  +1.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.122'
$ws.Range("E31").Value = '  +0.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1027'
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.710'
$ws.Range("E33").Value = '  +0.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.670'
$ws.Range("E34").Value = '  +1.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2298'
$ws.Range("E35").Value = '  +10.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06554'
$ws.Range("E36").Value = '  +8.86%  '
$ws.Range("E37").Value = '  +2.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.846'
$ws.Range("E38").Value = '  +2.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.071'
$ws.Range("E39").Value = '  +1.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.48'
$ws.Range("E40").Value = '  +0.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6298'
$ws.Range("E41").Value = '  +0.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.216'
$ws.Range("E42").Value = '  +6.09%  '
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.394'
$ws.Range("E44").Value = '  -1.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.47'
$ws.Range("E45").Value = '  +0.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5924'
$ws.Range("E46").Value = '  +1.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.669'
$ws.Range("E47").Value = '  +0.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '126.30'
$ws.Range("E48").Value = '  +3.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.979'
$ws.Range("E49").Value = '  +3.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.171'
$ws.Range("E50").Value = '  +3.59%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06927'
$ws.Range("E51").Value = '  +2.58%  '
